# Update the answer table: replace each populated cell's text with the
# newly generated three-digit ÷ one-digit division problem/answer.
#
# The source table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17
# (1-based) contain text. We address cells positionally via
# Table.Cell(row, col) so formatting (rFonts/sz) carried by the existing
# run is preserved automatically.

$d = $word.ActiveDocument
$t = $d.Tables(1)

# Row 1
$t.Cell(1, 1).Range.Text = "200÷6=33, 2"
$t.Cell(1, 2).Range.Text = "936÷6=156, 0"
$t.Cell(1, 3).Range.Text = "273÷7=39, 0"
$t.Cell(1, 4).Range.Text = "249÷8=31, 1"
$t.Cell(1, 5).Range.Text = "946÷8=118, 2"

# Row 5
$t.Cell(5, 1).Range.Text = "717÷8=89, 5"
$t.Cell(5, 2).Range.Text = "447÷6=74, 3"
$t.Cell(5, 3).Range.Text = "624÷8=78, 0"
$t.Cell(5, 4).Range.Text = "712÷2=356, 0"
$t.Cell(5, 5).Range.Text = "972÷5=194, 2"

# Row 9
$t.Cell(9, 1).Range.Text = "930÷5=186, 0"
$t.Cell(9, 2).Range.Text = "928÷3=309, 1"
$t.Cell(9, 3).Range.Text = "382÷7=54, 4"
$t.Cell(9, 4).Range.Text = "826÷3=275, 1"
$t.Cell(9, 5).Range.Text = "910÷7=130, 0"

# Row 13
$t.Cell(13, 1).Range.Text = "202÷5=40, 2"
$t.Cell(13, 2).Range.Text = "745÷2=372, 1"
$t.Cell(13, 3).Range.Text = "820÷7=117, 1"
$t.Cell(13, 4).Range.Text = "108÷2=54, 0"
$t.Cell(13, 5).Range.Text = "250÷2=125, 0"

# Row 17
$t.Cell(17, 1).Range.Text = "405÷8=50, 5"
$t.Cell(17, 2).Range.Text = "397÷7=56, 5"
$t.Cell(17, 3).Range.Text = "783÷4=195, 3"
$t.Cell(17, 4).Range.Text = "287÷3=95, 2"
$t.Cell(17, 5).Range.Text = "501÷6=83, 3"
